$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "fantasy points" values (currently in column E, rows
# 2-17) before they get overwritten, so they can be moved into new column G.
$oldFantasyPoints = @{}
for ($r = 2; $r -le 17; $r++) {
    $oldFantasyPoints[$r] = $ws.Cells.Item($r, 5).Value2
}

# The scraper re-ran and now also reports height/weight, inserted ahead of
# fantasy points - so "fantasy points" shifts out to the new column G.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Match the header style (bold font, thin border, centered/top aligned) used
# by the other header cells (B1:E1).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Populate the new height/weight columns with the scraped values and move the
# previous fantasy-points values into the new column G.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5
    $ws.Cells.Item($r, 6).Value = 267
    $ws.Cells.Item($r, 7).Value = $oldFantasyPoints[$r]
}
